# "major accuracy check update"
# - Shared string "E7420" -> "E7420L" (sample id in column G, rows 2-13)
# - Column H (accuracy-check flag) rows 2-13: was a volatile =FALSE() formula,
#   now stored as a plain boolean literal FALSE
# - Selection moves from G2:G13 to H2:H13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update every cell sharing the "E7420" string together so the workbook's
# shared-string table is edited in place instead of forking a new entry.
$ws.Range("G2:G13").Value = "E7420L"

# Replace the =FALSE() formula cells with literal boolean FALSE values.
$ws.Range("H2:H13").Value = $false

# Match the author's final selection (H2:H13, active cell H2).
$ws.Range("H2:H13").Select() | Out-Null
